# Rename the inline logo pictures in the document's headers/footers.
#
#   footer (first page)  : image2.png -> image1.png   (Pearson logo)
#   footer (primary)     : image2.png -> image1.png   (Pearson logo)
#   header (first page)  : image1.jpg -> image2.jpg   (BTEC logo)
#
# `InlineShape` has no writable `Name` property (matches real Word's
# object model), so the standard idiom is used: convert the inline
# picture to a floating `Shape`, rename it there, then convert it back
# to an inline picture - this preserves the wp:inline layout while
# updating the shape's name.

function Rename-InlineLogo($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# First-page footer (footer1.xml) - Pearson logo.
$firstPageFooterShape = $section.Footers.Item(2).Range.InlineShapes.Item(1)
Rename-InlineLogo $firstPageFooterShape "image1.png"

# Primary (other pages) footer (footer2.xml) - Pearson logo.
$primaryFooterShape = $section.Footers.Item(1).Range.InlineShapes.Item(1)
Rename-InlineLogo $primaryFooterShape "image1.png"

# First-page header (header1.xml) - BTEC logo.
$firstPageHeaderShape = $section.Headers.Item(2).Range.InlineShapes.Item(1)
Rename-InlineLogo $firstPageHeaderShape "image2.jpg"
